# Reorder the "Recorded By" (column G) values so that "System" is moved to
# the front of the comma-separated list, keeping the relative order of the
# remaining entries the same.
#
# e.g. "backup@backdoor.com, System, system" -> "System, backup@backdoor.com, system"
#      "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
# Rows whose value does not contain "System" as a distinct comma-separated
# entry (or that only contain a single entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }

    $parts = $value -split ', '

    if ($parts.Count -gt 1) {
        $systemIndex = -1
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($parts[$i].Equals('System')) {
                $systemIndex = $i
                break
            }
        }

        if ($systemIndex -ge 0) {
            $newParts = @('System')
            for ($i = 0; $i -lt $parts.Count; $i++) {
                if ($i -ne $systemIndex) {
                    $newParts += $parts[$i]
                }
            }
            $newValue = $newParts -join ', '
            if (-not $newValue.Equals($value)) {
                $cell.Value2 = $newValue
            }
        }
    }
}
